$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supervisor Name (row 6, merged G6:I6) - previously blank
$ws.Range("G6").Value = "Ankita Gangotra"

# Supervisor sign-off initials (row 27, merged A27:C27) - previously blank
$ws.Range("A27").Value = "A.G"

# Supervisor sign-off date (row 27, merged D27:E27) - previously blank
$ws.Range("D27").Value = (Get-Date -Year 2014 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("D27").NumberFormat = "mm-dd-yy"

# Update the active selection to reflect where the user left off editing
$ws.Range("G27").Select() | Out-Null
